$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Volume(1h)" (E) columns hold numeric-looking text values
# (e.g. "593.17", "0.0000149") that must stay plain text, not become Number cells.
# Temporarily force the whole D2:E51 block to Text format before writing the new
# values, then restore the original ("Normal") cell style once all writes are done.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.330.65"
$ws.Range("E2").Value = "  -2.85%  "

$ws.Range("D3").Value = "3.778.07"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "593.17"
$ws.Range("E5").Value = "  -1.21%  "

$ws.Range("D6").Value = "165.71"
$ws.Range("E6").Value = "  -3.12%  "

$ws.Range("D7").Value = "3.776.38"
$ws.Range("E7").Value = "  -0.72%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -1.72%  "

$ws.Range("E10").Value = "  -3.10%  "

$ws.Range("E11").Value = "  -2.39%  "

$ws.Range("E12").Value = "  -1.23%  "

$ws.Range("E13").Value = "  -4.15%  "

$ws.Range("D14").Value = "35.78"
$ws.Range("E14").Value = "  -3.09%  "

$ws.Range("D15").Value = "4.413.37"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "3.790.17"
$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("D17").Value = "67.348.57"
$ws.Range("E17").Value = "  -2.79%  "

$ws.Range("D18").Value = "17.93"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").Value = "6.95"
$ws.Range("E20").Value = "  -2.07%  "

$ws.Range("D21").Value = "10.21"
$ws.Range("E21").Value = "  -8.16%  "

$ws.Range("D22").Value = "456.07"
$ws.Range("E22").Value = "  -3.41%  "

$ws.Range("E23").Value = "  -1.61%  "

$ws.Range("D24").Value = "0.0000149"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").Value = "83.31"
$ws.Range("E25").Value = "  -1.89%  "

$ws.Range("D26").Value = "11.83"
$ws.Range("E26").Value = "  -3.31%  "

$ws.Range("E27").Value = "  -5.88%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").Value = "9.92"
$ws.Range("E29").Value = "  -3.64%  "

$ws.Range("E30").Value = "  -2.17%  "

$ws.Range("D31").Value = "29.74"
$ws.Range("E31").Value = "  -2.13%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "2.18"
$ws.Range("E32").Value = "  -3.09%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "7.17"
$ws.Range("E33").Value = "  -4.22%  "

$ws.Range("D34").Value = "9.14"
$ws.Range("E34").Value = "  -3.02%  "

$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("D36").Value = "3.733.08"

$ws.Range("D37").Value = "0.0993"
$ws.Range("E37").Value = "  -2.98%  "

$ws.Range("E38").Value = "  -2.04%  "

$ws.Range("D39").Value = "3.26"
$ws.Range("E39").Value = "  -7.50%  "

$ws.Range("E40").Value = "  -1.79%  "

$ws.Range("D41").Value = "5.70"
$ws.Range("E41").Value = "  -3.02%  "

$ws.Range("E42").Value = "  -0.15%  "

$ws.Range("D44").Value = "43.55"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("E45").Value = "  -4.47%  "

$ws.Range("D46").Value = "46.87"
$ws.Range("E46").Value = "  +1.60%  "

$ws.Range("D47").Value = "8.33"
$ws.Range("E47").Value = "  -3.96%  "

$ws.Range("D48").Value = "147.44"
$ws.Range("E48").Value = "  +1.85%  "

$ws.Range("D49").Value = "391.19"
$ws.Range("E49").Value = "  -2.99%  "

$ws.Range("D50").Value = "1.81"
$ws.Range("E50").Value = "  -8.27%  "

$ws.Range("D51").Value = "2.747.72"
$ws.Range("E51").Value = "  +1.76%  "

$dataRange.Style = "Normal"
